$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'280.11"
$ws.Range("E2").Value = "'5.39%"
$ws.Range("D3").Value = "'26.94"
$ws.Range("E3").Value = "'0.90%"
$ws.Range("D4").Value = "'4.946"
$ws.Range("E4").Value = "'5.12%"
$ws.Range("D5").Value = "'0.06417"
$ws.Range("E6").Value = "'4.77%"
$ws.Range("D7").Value = "'3.349"
$ws.Range("E7").Value = "'5.61%"
$ws.Range("D8").Value = "'0.8884"
$ws.Range("E8").Value = "'4.50%"
$ws.Range("D9").Value = "'1.027"
$ws.Range("E9").Value = "'13.46%"
$ws.Range("D10").Value = "'0.1498"
$ws.Range("E10").Value = "'6.39%"
$ws.Range("D11").Value = "'0.05197"
$ws.Range("E11").Value = "'7.19%"
$ws.Range("D12").Value = "'0.07405"
$ws.Range("E12").Value = "'4.27%"
$ws.Range("D13").Value = "'0.03137"
$ws.Range("E13").Value = "'-0.13%"
$ws.Range("D14").Value = "'0.09065"
$ws.Range("E14").Value = "'0.42%"
$ws.Range("D15").Value = "'0.001580"
$ws.Range("E15").Value = "'2.57%"
$ws.Range("D16").Value = "'0.0006313"
$ws.Range("E16").Value = "'4.24%"
$ws.Range("D17").Value = "'0.006058"
$ws.Range("E17").Value = "'1.80%"
$ws.Range("D18").Value = "'3.495"
$ws.Range("E18").Value = "'1.09%"
$ws.Range("D19").Value = "'2.297"
$ws.Range("E19").Value = "'0.87%"
$ws.Range("E20").Value = "'0.78%"
$ws.Range("E21").Value = "'2.47%"
$ws.Range("D22").Value = "'3.933"
$ws.Range("E22").Value = "'-4.20%"
$ws.Range("D23").Value = "'0.04370"
$ws.Range("E23").Value = "'3.08%"
$ws.Range("D24").Value = "'0.001181"
$ws.Range("E24").Value = "'-0.28%"
$ws.Range("D25").Value = "'0.003692"
$ws.Range("E25").Value = "'-10.74%"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'0.00%"
$ws.Range("D27").Value = "'0.0001696"
$ws.Range("E27").Value = "'0.83%"
$ws.Range("D40").Value = "'0.04112"
$ws.Range("E40").Value = "'4.76%"
$ws.Range("D41").Value = "'0.006641"
$ws.Range("E41").Value = "'59.24%"
$ws.Range("D42").Value = "'0.1179"
$ws.Range("E42").Value = "'5.79%"
$ws.Range("D43").Value = "'0.002361"
$ws.Range("E43").Value = "'11.87%"
$ws.Range("D44").Value = "'0.01256"
$ws.Range("E44").Value = "'9.48%"
$ws.Range("D45").Value = "'0.00005266"
$ws.Range("E45").Value = "'3.18%"
$ws.Range("E46").Value = "'-0.05%"
$ws.Range("E47").Value = "'1,348.53%"
$ws.Range("D48").Value = "'0.02242"
$ws.Range("E48").Value = "'-8.42%"
$ws.Range("E49").Value = "'-0.05%"
$ws.Range("E50").Value = "'-0.12%"
